# Auto-generated edit script: updates crypto price/volume table
# per commit "Updated cryptos list on Thu Dec  7 12:50:44 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.275.06"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "2.251.06"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.96"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.641"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.61"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.440"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0955"
$ws.Range("E10").Value = "  -7.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.48"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.96"
$ws.Range("E12").Value = "  +3.56%  "
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").Value = "2.583.54"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.97"
$ws.Range("E15").Value = "  -5.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.04"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.826"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "2.247.95"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "43.221.58"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("D20").Value = "0.0₃0964"
$ws.Range("E20").Value = "  -5.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.15"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.08"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.17"
$ws.Range("E23").Value = "  -3.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.77"
$ws.Range("E25").Value = "  +15.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.71"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.09"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.64"
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("E32").Value = "  -4.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.93"
$ws.Range("E34").Value = "  +3.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0678"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.91"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.62"
$ws.Range("E37").Value = "  -6.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.35"
$ws.Range("E38").Value = "  -6.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.27"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0249"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.73"
$ws.Range("E42").Value = "  +4.28%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.27"
$ws.Range("E43").Value = "  -1.15%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.48"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.87"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.25"
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.18"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0933"
$ws.Range("E48").Value = "  -2.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000207"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.432.51"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("E51").Value = "  -1.61%  "
